# "substract 19 to the stock series"
# - Subtract 19 from every value in the stock-series column (B2:B62) of the
#   "stock_driven" sheet. The originally-computed cells B11:B62 held a
#   formula (=previous+13); after the edit they hold plain literal values,
#   so we overwrite with .Value (not .Formula) to drop the formulas too.
# - Make "stock_driven" the active/selected sheet (it was "inflow_driven").
# - Update the view/selection remembered for each sheet to match where the
#   author ended up after editing.

$wb = $excel.ActiveWorkbook

$wsInflow = $wb.Worksheets.Item("inflow_driven")
$wsStock  = $wb.Worksheets.Item("stock_driven")

# --- Update the stock series values (subtract 19), replacing formulas with
#     plain numbers in the process. -----------------------------------------
$newValues = @(1,2,11,20,26,31,32,33,81,94,107,120,133,146,159,172,185,198,211,224,237,250,263,276,289,302,315,328,341,354,367,380,393,406,419,432,445,458,471,484,497,510,523,536,549,562,575,588,601,614,627,640,653,666,679,692,705,718,731,744,757)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 2 + $i
    $wsStock.Cells.Item($row, 2).Value = $newValues[$i]
}

# --- Move the "selected sheet" from inflow_driven to stock_driven ----------
$wsStock.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$wsStock.Range("F54").Select()

# --- Nudge the inflow chart's frame by a hair (matches the tiny rounding
#     drift seen after the source data changed). ---------------------------
$co = $wsInflow.ChartObjects().Item(1)
$emuNudge = 360.0 / 12700.0
$co.Width = $co.Width - $emuNudge
$co.Height = $co.Height - $emuNudge
